$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'60.147.51"
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Cells.Item(3, 4).Value = "'2.408.02"
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "'557.98"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "
$ws.Cells.Item(6, 4).Value = "'135.34"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Cells.Item(12, 4).Value = "'0.348"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Cells.Item(14, 4).Value = "'2.837.41"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Cells.Item(15, 4).Value = "'60.022.67"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Range("E15").Value = "  +0.27%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Cells.Item(17, 4).Value = "'2.455.85"
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Range("E17").Value = "  +2.56%  "
$ws.Cells.Item(18, 4).Value = "'11.21"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Range("E18").Value = "  -0.82%  "
$ws.Cells.Item(19, 4).Value = "'4.51"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Cells.Item(20, 4).Value = "'326.17"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Range("E20").Value = "  -1.28%  "
$ws.Cells.Item(21, 4).Value = "'6.79"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -2.13%  "
$ws.Cells.Item(24, 4).Value = "'0.172"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Cells.Item(25, 4).Value = "'8.54"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Range("E25").Value = "  -3.42%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Cells.Item(28, 4).Value = "'1.80"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Range("E28").Value = "  +1.33%  "
$ws.Range("E29").Value = "  -1.23%  "
$ws.Cells.Item(30, 4).Value = "'170.58"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +5.37%  "
$ws.Range("E33").Value = "  -1.73%  "
$ws.Cells.Item(34, 4).Value = "'18.40"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Range("E34").Value = "  -1.48%  "
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(35, 4).Value = "'0.999"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(36, 4).Value = "'1.33"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Range("E36").Value = "  +3.02%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  -0.84%  "
$ws.Cells.Item(39, 4).Value = "'324.64"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Range("E39").Value = "  +3.57%  "
$ws.Range("E40").Value = "  -0.65%  "
$ws.Cells.Item(41, 4).Value = "'38.55"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Cells.Item(42, 4).Value = "'148.52"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Range("E42").Value = "  +6.86%  "
$ws.Range("E43").Value = "  -3.32%  "
$ws.Range("E44").Value = "  -0.42%  "
$ws.Cells.Item(45, 4).Value = "'19.88"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Cells.Item(46, 4).Value = "'0.0515"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("E48").Value = "  -1.58%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("E50").Value = "  -1.29%  "
$ws.Range("E51").Value = "  -0.87%  "
